# "finished first draft of new layout"
# Remove the "DEFICIENCY" column from the "Deficiencies" sheet / Table1,
# leaving ID, ROOM, COMPONENT, NOTES - and make the Deficiencies sheet
# the active tab (with its NOTES column selected), matching the committed
# layout change.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Deficiencies")
$lo = $ws.ListObjects.Item("Table1")

# Drop the DEFICIENCY column (column D) and shift NOTES left into its place.
$ws.Columns.Item(4).Delete()

# Re-sync the table definition to the now-4-column range.
$lo.Resize($ws.Range("A1:D1048576"))

# The header cell text write keeps the ListObject's column name in sync
# with the sheet (NOTES is now the 4th/last column of Table1).
$ws.Range("D1").Value = "NOTES"

# "General Info" was the active tab before; it no longer is.
$wsGeneral = $wb.Worksheets.Item("General Info")

# Make "Deficiencies" the active sheet/tab, with the new NOTES column
# selected (mirrors where Excel leaves the cursor after a column delete).
$ws.Activate() | Out-Null
$ws.Range("D1:D1048576").Select() | Out-Null
